# Applies the cryptos.xlsx data refresh (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / link / name / percentage cells ---
# (safe to assign directly -- Excel keeps them as text since they are
#  not parseable as a single numeric literal)
$ws.Range("D2").Value = "25.998.34"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.641.93"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.671.39"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.870.26"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "26.091.38"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "1.128.79"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "1.779.53"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E51").Value = "  -0.41%  "

# --- Price cells whose new text *looks* numeric ---
# Excel auto-converts a plain numeric-looking string assigned to a
# General-formatted cell into a real number. The source data keeps
# these as plain text, so force Text format for the assignment, then
# restore the cell style to Normal/General so formatting is unaffected.
$numericPriceCells = @("D5", "D9", "D17", "D20", "D24", "D27", "D28", "D29", "D42", "D46", "D47", "D49", "D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "215.85"
$ws.Range("D9").Value = "0.0638"
$ws.Range("D17").Value = "63.37"
$ws.Range("D20").Value = "194.75"
$ws.Range("D24").Value = "1.79"
$ws.Range("D27").Value = "143.13"
$ws.Range("D28").Value = "6.87"
$ws.Range("D29").Value = "15.52"
$ws.Range("D42").Value = "99.17"
$ws.Range("D46").Value = "56.63"
$ws.Range("D47").Value = "0.0524"
$ws.Range("D49").Value = "7.75"
$ws.Range("D51").Value = "0.0954"
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
